# Apply update: insert new rows at final positions 17, 25, 26, 29, 34
# and rewrite rows 17-37 with the new dataset (post 'atualizacao post empresas com nova base').

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows, ascending by final row index so each insertion lands correctly
# relative to rows already shifted by earlier insertions.
$insertPositions = @(17, 25, 26, 29, 34)
foreach ($pos in $insertPositions) {
    $ws.Rows.Item($pos).Insert()
}

# Update the dimension-visible range A1:D37 with final values for rows 17-37.
$ws.Cells.Item(17, 1).Value = 900
$ws.Cells.Item(17, 2).Value = '71.473.820/0012-12'
$ws.Cells.Item(17, 3).Value = 'MILLENNIUM CCVM S/A, SUCESSORA DA GAMEX SECURITIES CCVM LTDA.'
$ws.Cells.Item(17, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2002&NumProc=15'

$ws.Cells.Item(18, 1).Value = 1456
$ws.Cells.Item(18, 2).Value = '47.894.290/0001-28'
$ws.Cells.Item(18, 3).Value = 'PLANIN AUDITORES INDEPENDENTES  S/C'
$ws.Cells.Item(18, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2001&NumProc=8739'

$ws.Cells.Item(19, 1).Value = 397
$ws.Cells.Item(19, 2).Value = '09.143.363/0001-50'
$ws.Cells.Item(19, 3).Value = 'EASE ESCRITÓRIO DE AUDITORIA INDEPENDENTE S/C'
$ws.Cells.Item(19, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2013&NumProc=13481'

$ws.Cells.Item(20, 1).Value = 411
$ws.Cells.Item(20, 2).Value = '62.030.762/0001-98'
$ws.Cells.Item(20, 3).Value = 'AKW AUDITORES INDEPENDENTES S/S'
$ws.Cells.Item(20, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2013&NumProc=4362'

$ws.Cells.Item(21, 1).Value = 1165
$ws.Cells.Item(21, 2).Value = '67.634.717/0001-66'
$ws.Cells.Item(21, 3).Value = 'BWEL AUDITORES INDEPENDENTES SOCIEDADE SIMPLES'
$ws.Cells.Item(21, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2017&NumProc=2253'

$ws.Cells.Item(22, 1).Value = 427
$ws.Cells.Item(22, 2).Value = '11.245.719/0003-70'
$ws.Cells.Item(22, 3).Value = 'DIRECTA AUDITORES'
$ws.Cells.Item(22, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2013&NumProc=5682'

$ws.Cells.Item(23, 1).Value = 1688
$ws.Cells.Item(23, 2).Value = '58.214.958/0001-65'
$ws.Cells.Item(23, 3).Value = 'SOC TEC AUDITORIA SOMATEC SC'
$ws.Cells.Item(23, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2000&NumProc=6'

$ws.Cells.Item(24, 1).Value = 844
$ws.Cells.Item(24, 2).Value = '31.622.483/0001-90'
$ws.Cells.Item(24, 3).Value = 'GUILHERME FONTES FILMES LTDA.'
$ws.Cells.Item(24, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2000&NumProc=12'

$ws.Cells.Item(25, 1).Value = 1035
$ws.Cells.Item(25, 2).Value = '00.469.585/0001-93'
$ws.Cells.Item(25, 3).Value = 'FACEB - FUNDAÇÃO DE ASSISTÊNCIA DOS EMPREGADOS DA CEB'
$ws.Cells.Item(25, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=1999&NumProc=28'

$ws.Cells.Item(26, 1).Value = 1035
$ws.Cells.Item(26, 2).Value = '17.393.471/0001-13'
$ws.Cells.Item(26, 3).Value = 'PRATA DTVM LTDA. (atual Prata Consultoria e Assessoria Empresarial Ltda.)'
$ws.Cells.Item(26, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=1999&NumProc=28'

$ws.Cells.Item(27, 1).Value = 52
$ws.Cells.Item(27, 2).Value = '04.612.682/0001-44'
$ws.Cells.Item(27, 3).Value = 'INTERTRADING AGRONEGÓCIOS LTDA.'
$ws.Cells.Item(27, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2007&NumProc=4414'

$ws.Cells.Item(28, 1).Value = 21
$ws.Cells.Item(28, 2).Value = '27.901.719/0001-50'
$ws.Cells.Item(28, 3).Value = 'INSTITUTO AERUS DE SEGURIDADE SOCIAL'
$ws.Cells.Item(28, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2007&NumProc=1176'

$ws.Cells.Item(29, 1).Value = 1563
$ws.Cells.Item(29, 2).Value = '43.214.485/0001-29'
$ws.Cells.Item(29, 3).Value = 'SOCIVAL AUDITORIA INDEP SC'
$ws.Cells.Item(29, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2004&NumProc=7001'

$ws.Cells.Item(30, 1).Value = 555
$ws.Cells.Item(30, 2).Value = '64.920.416/0001-00'
$ws.Cells.Item(30, 3).Value = 'NORMAS AUDITORES INDEPENDENTES'
$ws.Cells.Item(30, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2015&NumProc=11941'

$ws.Cells.Item(31, 1).Value = 1565
$ws.Cells.Item(31, 2).Value = '43.729.789/0001-29'
$ws.Cells.Item(31, 3).Value = 'PERMALI DO BRASIL INDÚSTRIA E COMÉRCIO LTDA'
$ws.Cells.Item(31, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2005&NumProc=33'

$ws.Cells.Item(32, 1).Value = 1556
$ws.Cells.Item(32, 2).Value = '04.565.230/0002-30'
$ws.Cells.Item(32, 3).Value = 'I.B. Sabbá S/A'
$ws.Cells.Item(32, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2004&NumProc=4627'

$ws.Cells.Item(33, 1).Value = 217
$ws.Cells.Item(33, 2).Value = '05.723.617/0001-59'
$ws.Cells.Item(33, 3).Value = 'MAPFRE DTVM S.A.'
$ws.Cells.Item(33, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2010&NumProc=17292'

$ws.Cells.Item(34, 1).Value = 1710
$ws.Cells.Item(34, 2).Value = '14.629.882/0001-63'
$ws.Cells.Item(34, 3).Value = 'CAPITAL ASSESSORIA FINANCEIRA LTDA. (ATUAL CAPITAL ASSESSORIA E EMPREENDIMENTOS LTDA.)'
$ws.Cells.Item(34, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2002&NumProc=6'

$ws.Cells.Item(35, 1).Value = 642
$ws.Cells.Item(35, 2).Value = '05.706.592/0001-85'
$ws.Cells.Item(35, 3).Value = 'BANCO BOZANO, SIMONSEN S/A'
$ws.Cells.Item(35, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=1999&NumProc=10'

$ws.Cells.Item(36, 1).Value = 2
$ws.Cells.Item(36, 2).Value = '00.659.559/0002-09'
$ws.Cells.Item(36, 3).Value = 'MASTER CORRETORA DE MERCADORIAS LTDA.'
$ws.Cells.Item(36, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=IA&Ano=2006&NumProc=1'

$ws.Cells.Item(37, 1).Value = 4
$ws.Cells.Item(37, 2).Value = '00.016.087/6747-72'
$ws.Cells.Item(37, 3).Value = 'MARCOS LEVY'
$ws.Cells.Item(37, 4).Value = 'http://sistemas.cvm.gov.br/asp/cvmwww/inqueritos/DetPasAndamento.asp?sg_uf=RJ&Ano=2006&NumProc=8625'
